$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-12 with the new year labels and values
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 24872.5952
$ws.Range("C2").Value = 47003.8686
$ws.Range("D2").Value = 50349.1708

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 30264.78871
$ws.Range("C3").Value = 58833.66012
$ws.Range("D3").Value = 63595.06712

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 36123.45204
$ws.Range("C4").Value = 68539.05628
$ws.Range("D4").Value = 75529.71696000001

$ws.Range("A5").Value = "2014年"
$ws.Range("B5").Value = 47913.33985
$ws.Range("C5").Value = 90525.66266
$ws.Range("D5").Value = 100459.18004
$ws.Range("E5").Value = 7011.35495

$ws.Range("A6").Value = "2015年"
$ws.Range("B6").Value = 54195.92887
$ws.Range("C6").Value = 98744.24892
$ws.Range("D6").Value = 110114.37228
$ws.Range("E6").Value = 7007.55063

$ws.Range("A7").Value = "2016年"
$ws.Range("B7").Value = 60946.43791
$ws.Range("C7").Value = 109658.47806
$ws.Range("D7").Value = 121505.27937
$ws.Range("E7").Value = 11846.80131

$ws.Range("A8").Value = "2017年"
$ws.Range("B8").Value = 67726.45039
$ws.Range("C8").Value = 122841.24935
$ws.Range("D8").Value = 136900.6907
$ws.Range("E8").Value = 14059.44135

$ws.Range("A9").Value = "2018年"
$ws.Range("B9").Value = 74761.86724000001
$ws.Range("C9").Value = 141757.74955
$ws.Range("D9").Value = 159218.09563
$ws.Range("E9").Value = 11647.17383

$ws.Range("A10").Value = "2019年"
$ws.Range("B10").Value = 80154.78904
$ws.Range("C10").Value = 159024.95319
$ws.Range("D10").Value = 176472.89486
$ws.Range("E10").Value = 13190.21102

$ws.Range("A11").Value = "2020年"
$ws.Range("B11").Value = 85100.23762
$ws.Range("C11").Value = 178466.9915
$ws.Range("D11").Value = 197934.67822
$ws.Range("E11").Value = 15482.42515

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 90528.51068000001
$ws.Range("C12").Value = 200126.51499
$ws.Range("D12").Value = 220767.25259
$ws.Range("E12").Value = 16294.52733

# Remove the now-unused trailing rows (13-20) so the sheet shrinks to A1:E12
$ws.Range("A13:E20").EntireRow.Delete()
